{"js": "// Consolidate the split (one-word-per-run) text of the Title, Author and\n// Date paragraphs at the top of the document into a single run each,\n// matching the canonical OOXML produced by the later Bookdown build.\n// (The paragraph's visible text is unchanged - only the run structure is\n// collapsed from \"one run per word/space\" down to a single run.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Style name -> the (unchanged) full text that should end up in one run.\nconst targets = {\n  \"Title\": \"Test 003: Repeat an environment with the correct numbering\",\n  \"Author\": \"Emma Cliffe, Skills Centre: MASH, University of Bath\",\n  \"Date\": \"August 2020\",\n};\n\nconst remaining = new Set(Object.keys(targets));\nfor (const paragraph of paragraphs.items) {\n  if (remaining.size === 0) break;\n  const style = paragraph.style;\n  if (remaining.has(style)) {\n    paragraph.insertText(targets[style], \"Replace\");\n    remaining.delete(style);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Consolidate the split (one-word-per-run) text of the Title, Author and\n# Date paragraphs at the top of the document into a single run each,\n# matching the canonical OOXML produced by the later Bookdown build.\n# (The paragraph's visible text is unchanged - only the run structure is\n# collapsed from \"one run per word/space\" down to a single run.)\n\n$d = $word.ActiveDocument\n\n$targets = @{\n    \"Title\"  = \"Test 003: Repeat an environment with the correct numbering\"\n    \"Author\" = \"Emma Cliffe, Skills Centre: MASH, University of Bath\"\n    \"Date\"   = \"August 2020\"\n}\n\n$remainingCount = $targets.Count\n\nforeach ($p in $d.Paragraphs) {\n    if ($remainingCount -eq 0) { break }\n    $styleName = $p.Style.NameLocal\n    if ($targets.ContainsKey($styleName)) {\n        $r = $d.Range($p.Range.Start, $p.Range.End)\n        $r.Text = $targets[$styleName]\n        $targets.Remove($styleName)\n        $remainingCount = $remainingCount - 1\n    }\n}\n"}
